$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be read as numbers,
# so they stay text like the rest of the Price column.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.920.58"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "2.231.91"
$ws.Range("E3").Value = "  -2.07%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").Value = "311.97"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("D6").Value = "98.50"
$ws.Range("E6").Value = "  -4.94%  "
$ws.Range("D7").Value = "0.568"
$ws.Range("E7").Value = "  -3.62%  "
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -7.33%  "
$ws.Range("D10").Value = "36.06"
$ws.Range("E10").Value = "  -5.70%  "
$ws.Range("E11").Value = "  -3.24%  "
$ws.Range("D12").Value = "7.33"
$ws.Range("E12").Value = "  -6.96%  "
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("D14").Value = "2.577.17"
$ws.Range("E14").Value = "  -1.44%  "
$ws.Range("D15").Value = "2.240.19"
$ws.Range("E15").Value = "  -1.39%  "
$ws.Range("D16").Value = "0.834"
$ws.Range("E16").Value = "  -5.21%  "
$ws.Range("D17").Value = "14.07"
$ws.Range("E17").Value = "  -3.73%  "
$ws.Range("D18").Value = "43.819.90"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  -10.36%  "
$ws.Range("D20").Value = "0.0₃0959"
$ws.Range("E20").Value = "  -4.08%  "
$ws.Range("E21").Value = "  -5.84%  "
$ws.Range("D22").Value = "64.77"
$ws.Range("E22").Value = "  -2.25%  "
$ws.Range("D23").Value = "2.99"
$ws.Range("E23").Value = "  -7.45%  "
$ws.Range("D24").Value = "232.39"
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("E25").Value = "  -9.56%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "10.13"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").Value = "36.77"
$ws.Range("E29").Value = "  -7.31%  "
$ws.Range("D30").Value = "5.91"
$ws.Range("E30").Value = "  -9.47%  "
$ws.Range("D31").Value = "157.52"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").Value = "19.85"
$ws.Range("E32").Value = "  -3.65%  "
$ws.Range("D33").Value = "0.0827"
$ws.Range("E33").Value = "  -6.87%  "
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "3.19"
$ws.Range("E35").Value = "  -7.04%  "
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("E37").Value = "  -7.84%  "
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("D39").Value = "15.62"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("E40").Value = "  -8.97%  "
$ws.Range("D41").Value = "4.05"
$ws.Range("E41").Value = "  -10.63%  "
$ws.Range("E42").Value = "  -7.31%  "
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "1.712.51"
$ws.Range("E44").Value = "  -5.89%  "
$ws.Range("D45").Value = "0.193"
$ws.Range("E45").Value = "  -8.10%  "
$ws.Range("D46").Value = "79.91"
$ws.Range("E46").Value = "  -7.66%  "
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("E48").Value = "  -6.50%  "
$ws.Range("D49").Value = "72.59"
$ws.Range("E49").Value = "  -5.81%  "
$ws.Range("D50").Value = "101.28"
$ws.Range("E50").Value = "  -3.32%  "
$ws.Range("D51").Value = "56.12"
$ws.Range("E51").Value = "  -6.55%  "
